$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.143.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.815.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  +0.77%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5912'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.007'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("E8").Value = '  -4.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06810'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.96'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07542'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.660'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6174'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009549'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '75.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.932.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.429'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -9.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.739'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.007'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.797'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1269'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06370'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.411'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.438'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.737'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.706'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.39%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.694'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.92%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.070'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.03%  '
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6330'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.753'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01719'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.87%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.459'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.125.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8785'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.007'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.972.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000115'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.32%  '
$ws.Range("E47").Value = '  -3.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05503'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("E49").Value = '  -0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.322'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3571'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.82%  '
